$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the diff (Price / Volume(1h) columns, plus
# a couple of row swaps for Coin name / Link columns).

$ws.Range("D2").Value = "65.642.27"
$ws.Range("D3").Value = "2.650.08"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.626"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000197"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D15").Value = "3.127.18"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "65.467.44"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "2.640.82"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "155.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0605"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.634"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0993"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "0.0₆0250"
$ws.Range("E50").Value = "  +7.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.84%  "
